$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, preserving its original
# "General" style (no visible numeric coercion of digit-and-dot strings,
# no residual cell-level number-format override left behind).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "60.395.00"
Set-TextValue "E2" "  +6.32%  "
Set-TextValue "D3" "3.316.08"
Set-TextValue "E3" "  +2.01%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "408.89"
Set-TextValue "E5" "  +3.23%  "
Set-TextValue "D6" "113.75"
Set-TextValue "E6" "  +4.53%  "
Set-TextValue "D7" "3.309.22"
Set-TextValue "E7" "  +1.89%  "
Set-TextValue "E8" "  -1.47%  "
Set-TextValue "D9" "1.00"
Set-TextValue "E9" "  +0.06%  "
Set-TextValue "D10" "0.623"
Set-TextValue "E10" "  -0.35%  "
Set-TextValue "D11" "0.115"
Set-TextValue "E11" "  +18.14%  "
Set-TextValue "D12" "39.39"
Set-TextValue "E12" "  +0.26%  "
Set-TextValue "E13" "  -0.86%  "
Set-TextValue "D14" "3.834.09"
Set-TextValue "E14" "  +1.84%  "
Set-TextValue "D15" "8.16"
Set-TextValue "E15" "  -2.12%  "
Set-TextValue "D16" "19.03"
Set-TextValue "E16" "  -0.50%  "
Set-TextValue "D17" "3.310.60"
Set-TextValue "E17" "  +1.96%  "
Set-TextValue "D18" "60.291.24"
Set-TextValue "E18" "  +6.25%  "
Set-TextValue "D19" "0.997"
Set-TextValue "E19" "  -3.47%  "
Set-TextValue "D20" "10.78"
Set-TextValue "E20" "  +0.77%  "
Set-TextValue "D21" "0.0000115"
Set-TextValue "E21" "  +5.50%  "
Set-TextValue "D22" "3.34"
Set-TextValue "E22" "  -0.66%  "
Set-TextValue "D23" "12.36"
Set-TextValue "E23" "  -4.26%  "
Set-TextValue "D24" "295.95"
Set-TextValue "D25" "73.74"
Set-TextValue "E25" "  -0.61%  "
Set-TextValue "E26" "  -2.14%  "
Set-TextValue "D27" "29.12"
Set-TextValue "E27" "  +3.34%  "
Set-TextValue "D28" "7.57"
Set-TextValue "E28" "  +4.19%  "
Set-TextValue "D29" "4.24"
Set-TextValue "E29" "  -2.37%  "
Set-TextValue "E30" "  +2.18%  "
Set-TextValue "D31" "7.46"
Set-TextValue "E31" "  -2.52%  "
Set-TextValue "E32" "  +4.28%  "
Set-TextValue "E33" "  +0.05%  "
Set-TextValue "D34" "11.21"
Set-TextValue "E34" "  -0.70%  "
Set-TextValue "D35" "2.48"
Set-TextValue "E35" "  +16.35%  "
Set-TextValue "D36" "40.49"
Set-TextValue "E36" "  +1.83%  "
Set-TextValue "D37" "0.0489"
Set-TextValue "E37" "  +0.82%  "
Set-TextValue "D38" "52.07"
Set-TextValue "E38" "  +1.39%  "
Set-TextValue "D39" "0.998"
Set-TextValue "E39" "  -0.06%  "
Set-TextValue "E40" "  +5.86%  "
Set-TextValue "D41" "3.34"
Set-TextValue "E41" "  -3.95%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D42" "134.18"
Set-TextValue "E42" "  -1.85%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D43" "0.292"
Set-TextValue "E43" "  +3.89%  "
Set-TextValue "E44" "  -1.59%  "
Set-TextValue "E45" "  -0.87%  "
Set-TextValue "D46" "3.78"
Set-TextValue "E46" "  -4.70%  "
Set-TextValue "D47" "16.22"
Set-TextValue "E47" "  -4.71%  "
Set-TextValue "D48" "2.19"
Set-TextValue "E48" "  +2.59%  "
Set-TextValue "D49" "20.97"
Set-TextValue "E49" "  -5.96%  "
Set-TextValue "D50" "2.132.70"
Set-TextValue "E50" "  -1.16%  "
Set-TextValue "D51" "3.632.32"
Set-TextValue "E51" "  +1.63%  "
